# Update event attribute info table: reshuffle Attribute/Type rows (2-21),
# row 11 ("case:concept:name"/"str") is unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("case", "str"),
    @("human_workstation_green_button_pressed", "float"),
    @("lifecycle:state", "str"),
    @("operation_end_time", "datetime"),
    @("complete_service_time", "str"),
    @("identifier:id", "str"),
    @("process_model_id", "str"),
    @("concept:name", "str"),
    @("requested_service_url", "str"),
    @("case:concept:name", "str"),
    @("unsatisfied_condition_description", "str"),
    @("lifecycle:transition", "str"),
    @("planned_operation_time", "str"),
    @("event_id", "str"),
    @("SubProcessID", "str"),
    @("org:resource", "str"),
    @("current_task", "str"),
    @("time:timestamp", "datetime"),
    @("parameters", "dict"),
    @("response_status_code", "float")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
